# Apply the "msdescription" TEI style tweaks:
#   1. Normal style gains a <w:spacing w:after="120"/> paragraph property.
#   2. Two new custom paragraph styles - tei_collation and tei_extent -
#      are added to the style sheet, both based on Heading 4 with Normal
#      as the style for the following paragraph.

$d = $word.ActiveDocument

# --- 1. Normal style: add spacing-after ------------------------------------
$normal = $d.Styles.Item("Normal")
$normal.ParagraphFormat.SpaceAfter = 6

# --- 2. New custom styles: tei_collation, tei_extent ------------------------
$wdStyleTypeParagraph = 1

$collation = $d.Styles.Add("teicollation", $wdStyleTypeParagraph)
$collation.NameLocal = "tei_collation"
$collation.BaseStyle = "Heading4"
$collation.NextParagraphStyle = "Normal"
$collation.QuickStyle = $true

$extent = $d.Styles.Add("teiextent", $wdStyleTypeParagraph)
$extent.NameLocal = "tei_extent"
$extent.BaseStyle = "Heading4"
$extent.NextParagraphStyle = "Normal"
$extent.QuickStyle = $true
